$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 3.5
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 2.3
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.75
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("Z2").Value = 12
$ws.Range("AC2").Value = 7.5
$ws.Range("AV2").Value = 81
